$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing aspect/archetype strings ---
$ws.Range("J1").Value = "Stealthy"
$ws.Range("J5").Value = "Dropout Mage"

# --- Rework row 10 / add row 11 (new archetype block) ---
# Old state: only L10 = "Distant Doom"
# New state: row10 becomes "Sniper" archetype row, row11 becomes "Warlock" archetype row

$ws.Range("L10").ClearContents()

$ws.Range("A10").Value = "Sniper"
$ws.Range("C10").Value = "ALL"
$ws.Range("J10").Value = "ALL"
$ws.Range("O10").Value = "Fusilier"
$ws.Range("P10").Value = "ALL"
$ws.Range("R10").Value = "Wild Hunter"
$ws.Range("T10").Value = "Pact Archer"

$ws.Range("A11").Value = "Warlock"
$ws.Range("B11").Value = "Duskblade"
$ws.Range("E11").Value = "ALL"
$ws.Range("F11").Value = "Necromancer"
$ws.Range("H11").Value = "Warp Weaver"
$ws.Range("L11").Value = "ALL"
$ws.Range("T11").Value = "ALL"

# --- Selection moves to F12 in the saved view ---
$ws.Range("F12").Select()
